$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 139
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = $null

$ws.Range("H17").Value = 2473
$ws.Range("J17").Value = 2473
$ws.Range("L17").Value = 7419
$ws.Range("N17").Value = -7755

$ws.Range("I18").Value = 5818
$ws.Range("K18").Value = 5818
$ws.Range("M18").Value = -5534

$ws.Range("H43").Value = 4319.227
$ws.Range("I43").Value = 5961
$ws.Range("K43").Value = 5961
$ws.Range("M43").Value = -5892

$ws.Range("H98").Value = 4993.0347
$ws.Range("I98").Value = 5108.5
$ws.Range("K98").Value = 5108.5
$ws.Range("M98").Value = -3610.5

$ws.Range("H106").Value = 2744.8823
$ws.Range("I106").Value = 2668.7856
$ws.Range("K106").Value = 2668.7856
$ws.Range("M106").Value = -2037.7856

$ws.Range("H111").Value = 626.6
$ws.Range("I111").Value = 460.57144
$ws.Range("K111").Value = 1381.71432
$ws.Range("M111").Value = 1685.28568

$ws.Range("H113").Value = 1264
$ws.Range("I113").Value = 1149
$ws.Range("J113").Value = 1494
$ws.Range("K113").Value = 1149
$ws.Range("L113").Value = 1494
$ws.Range("M113").Value = 2105
$ws.Range("N113").Value = -8002

$ws.Range("H116").Value = 2325532
$ws.Range("I116").Value = 2325532
$ws.Range("K116").Value = 2325532
$ws.Range("M116").Value = -2322090

$ws.Range("H122").Value = 4993.0347
$ws.Range("I122").Value = 5108.5
$ws.Range("K122").Value = 15325.5
$ws.Range("M122").Value = -12875.5

$ws.Range("H127").Value = 1273.2
$ws.Range("I127").Value = 1091.75
$ws.Range("J127").Value = 1999
$ws.Range("K127").Value = 3275.25
$ws.Range("L127").Value = 5997
$ws.Range("M127").Value = 1684.75
$ws.Range("N127").Value = -15917

$ws.Range("H132").Value = 5962.759
$ws.Range("I132").Value = 5962.759
$ws.Range("K132").Value = 17888.277
$ws.Range("M132").Value = -15358.277

$ws.Range("H138").Value = 301416.8
$ws.Range("J138").Value = 484654.47
$ws.Range("L138").Value = 1453963.41
$ws.Range("N138").Value = -1464243.41

$ws.Range("H141").Value = 3352.2
$ws.Range("I141").Value = 2052.5557
$ws.Range("J141").Value = 5301.6665
$ws.Range("K141").Value = 6157.6671
$ws.Range("L141").Value = 15904.9995
$ws.Range("M141").Value = -977.6670999999997
$ws.Range("N141").Value = -26264.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 707.8889
$ws.Range("I2").Value = 666.3043
$ws.Range("J2").Value = 947
$ws.Range("K2").Value = 666.3043
$ws.Range("L2").Value = 947
$ws.Range("M2").Value = -553.3043
$ws.Range("N2").Value = -1173

$ws.Range("H32").Value = 3677.5244
$ws.Range("I32").Value = 2838.6753
$ws.Range("K32").Value = 2838.6753
$ws.Range("M32").Value = -2551.6753

$ws.Range("H45").Value = 28712.5
$ws.Range("I45").Value = 33685.383
$ws.Range("K45").Value = 33685.383
$ws.Range("M45").Value = -33308.383

$ws.Range("H61").Value = 6250
$ws.Range("I61").Value = 3139.1667
$ws.Range("K61").Value = 3139.1667
$ws.Range("M61").Value = -2927.1667

$ws.Range("H74").Value = 271504.47
$ws.Range("I74").Value = 618665
$ws.Range("J74").Value = 11134.083
$ws.Range("K74").Value = 618665
$ws.Range("L74").Value = 11134.083
$ws.Range("M74").Value = -617791
$ws.Range("N74").Value = -12882.083

$ws.Range("H77").Value = 271504.47
$ws.Range("I77").Value = 618665
$ws.Range("J77").Value = 11134.083
$ws.Range("K77").Value = 3093325
$ws.Range("L77").Value = 55670.415
$ws.Range("M77").Value = -3088957
$ws.Range("N77").Value = -64406.415

$ws.Range("H97").Value = 976.375
$ws.Range("I97").Value = 984.2381
$ws.Range("J97").Value = 921.3333
$ws.Range("K97").Value = 984.2381
$ws.Range("L97").Value = 921.3333
$ws.Range("M97").Value = -488.2381
$ws.Range("N97").Value = -1913.3333

$ws.Range("H116").Value = 707.8889
$ws.Range("I116").Value = 666.3043
$ws.Range("J116").Value = 947
$ws.Range("K116").Value = 666.3043
$ws.Range("L116").Value = 947
$ws.Range("M116").Value = 1627.6957
$ws.Range("N116").Value = -5535

$ws.Range("H132").Value = 2315.7576
$ws.Range("I132").Value = 1608.7391
$ws.Range("J132").Value = 3941.9
$ws.Range("K132").Value = 4826.2173
$ws.Range("L132").Value = 11825.7
$ws.Range("M132").Value = -2296.2173
$ws.Range("N132").Value = -16885.7

$ws.Range("H133").Value = 36430.25
$ws.Range("J133").Value = 36430.25
$ws.Range("L133").Value = 36430.25
$ws.Range("N133").Value = -41490.25

$ws.Range("H135").Value = 104737.75
$ws.Range("J135").Value = 104737.75
$ws.Range("L135").Value = 104737.75
$ws.Range("N135").Value = -114877.75

$ws.Range("H136").Value = 6250
$ws.Range("I136").Value = 3139.1667
$ws.Range("K136").Value = 9417.500100000001
$ws.Range("M136").Value = -6867.500100000001

$ws.Range("H139").Value = 106994.5
$ws.Range("J139").Value = 106994.5
$ws.Range("L139").Value = 106994.5
$ws.Range("N139").Value = -117274.5

$ws.Range("H140").Value = 76997.5
$ws.Range("J140").Value = 76997.5
$ws.Range("L140").Value = 76997.5
$ws.Range("N140").Value = -87357.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 707.8889
$ws.Range("I3").Value = 666.3043
$ws.Range("J3").Value = 947
$ws.Range("K3").Value = 666.3043
$ws.Range("L3").Value = 947
$ws.Range("M3").Value = -552.3043
$ws.Range("N3").Value = -1175

$ws.Range("H82").Value = 58717.91
$ws.Range("I82").Value = 40807.145
$ws.Range("K82").Value = 40807.145
$ws.Range("M82").Value = -40424.145

$ws.Range("H85").Value = 58717.91
$ws.Range("I85").Value = 40807.145
$ws.Range("K85").Value = 40807.145
$ws.Range("M85").Value = -39481.145

$ws.Range("H86").Value = 3465.5151
$ws.Range("I86").Value = 3226.16
$ws.Range("K86").Value = 3226.16
$ws.Range("M86").Value = -2103.16

$ws.Range("H89").Value = 3465.5151
$ws.Range("I89").Value = 3226.16
$ws.Range("K89").Value = 16130.8
$ws.Range("M89").Value = -10514.8

$ws.Range("H94").Value = 58824150
$ws.Range("I94").Value = 58824150
$ws.Range("K94").Value = 58824150
$ws.Range("M94").Value = -58823699

$ws.Range("H99").Value = 127384.875
$ws.Range("I99").Value = 167346.5
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 167346.5
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -165848.5
$ws.Range("N99").Value = -10496

$ws.Range("H105").Value = 10002498
$ws.Range("I105").Value = 590442
$ws.Range("K105").Value = 590442
$ws.Range("M105").Value = -588695

$ws.Range("H107").Value = 2959969.8
$ws.Range("I107").Value = 3497749.2
$ws.Range("K107").Value = 3497749.2
$ws.Range("M107").Value = -3495829.2

$ws.Range("H134").Value = 3558.6667
$ws.Range("I134").Value = 2817.2354
$ws.Range("J134").Value = 5359.2856
$ws.Range("K134").Value = 8451.706200000001
$ws.Range("L134").Value = 16077.8568
$ws.Range("M134").Value = -5916.706200000001
$ws.Range("N134").Value = -21147.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66668024
$ws.Range("I7").Value = 1538.4
$ws.Range("K7").Value = 1538.4
$ws.Range("M7").Value = -1425.4

$ws.Range("H16").Value = 1767.8182
$ws.Range("I16").Value = 1894.6
$ws.Range("K16").Value = 1894.6
$ws.Range("M16").Value = -1607.6

$ws.Range("H22").Value = 1286.1177
$ws.Range("I22").Value = 1019.4
$ws.Range("J22").Value = 1667.1428
$ws.Range("K22").Value = 1019.4
$ws.Range("L22").Value = 1667.1428
$ws.Range("M22").Value = -669.4
$ws.Range("N22").Value = -2367.1428

$ws.Range("H31").Value = 4277.6206
$ws.Range("I31").Value = 4026.7368
$ws.Range("J31").Value = 4399.846
$ws.Range("K31").Value = 4026.7368
$ws.Range("L31").Value = 4399.846
$ws.Range("M31").Value = -3731.7368
$ws.Range("N31").Value = -4989.846

$ws.Range("H34").Value = 4277.6206
$ws.Range("I34").Value = 4026.7368
$ws.Range("J34").Value = 4399.846
$ws.Range("K34").Value = 4026.7368
$ws.Range("L34").Value = 4399.846
$ws.Range("M34").Value = -3824.7368
$ws.Range("N34").Value = -4803.846

$ws.Range("H105").Value = 2563.182
$ws.Range("I105").Value = 1166.6666
$ws.Range("K105").Value = 1166.6666
$ws.Range("M105").Value = 580.3334

$ws.Range("H107").Value = 2632381.5
$ws.Range("I107").Value = 3334180
$ws.Range("J107").Value = 637.25
$ws.Range("K107").Value = 3334180
$ws.Range("L107").Value = 637.25
$ws.Range("M107").Value = -3332260
$ws.Range("N107").Value = -4477.25

$ws.Range("H113").Value = 1767.8182
$ws.Range("I113").Value = 1894.6
$ws.Range("K113").Value = 1894.6
$ws.Range("M113").Value = 275.4000000000001

$ws.Range("H122").Value = 3711.5
$ws.Range("I122").Value = 3226.077
$ws.Range("J122").Value = 4973.6
$ws.Range("K122").Value = 9678.231
$ws.Range("L122").Value = 14920.8
$ws.Range("M122").Value = -7228.231
$ws.Range("N122").Value = -19820.8

$ws.Range("H132").Value = 12824023
$ws.Range("I132").Value = 15628153
$ws.Range("K132").Value = 46884459
$ws.Range("M132").Value = -46881929

$ws.Range("H134").Value = 2974.3914
$ws.Range("I134").Value = 2417.9443
$ws.Range("K134").Value = 7253.8329
$ws.Range("M134").Value = -4718.8329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2789.8333
$ws.Range("I5").Value = 576.375
$ws.Range("J5").Value = 3594.7273
$ws.Range("K5").Value = 1729.125
$ws.Range("L5").Value = 10784.1819
$ws.Range("M5").Value = -1617.125
$ws.Range("N5").Value = -11008.1819

$ws.Range("H11").Value = 200209.5
$ws.Range("I11").Value = 125199.375
$ws.Range("J11").Value = 500250
$ws.Range("K11").Value = 375598.125
$ws.Range("L11").Value = 1500750
$ws.Range("M11").Value = -375458.125
$ws.Range("N11").Value = -1501030

$ws.Range("H56").Value = 7155.615
$ws.Range("I56").Value = 7155.615
$ws.Range("K56").Value = 7155.615
$ws.Range("M56").Value = -6625.615

$ws.Range("H70").Value = 1527
$ws.Range("I70").Value = 1527
$ws.Range("K70").Value = 4581
$ws.Range("M70").Value = -4266

$ws.Range("H73").Value = 1527
$ws.Range("I73").Value = 1527
$ws.Range("K73").Value = 4581
$ws.Range("M73").Value = -3489

$ws.Range("H108").Value = 1606.6666
$ws.Range("J108").Value = 4000
$ws.Range("L108").Value = 12000
$ws.Range("N108").Value = -17760

$ws.Range("H110").Value = 3635
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws.Range("H112").Value = 3048.75
$ws.Range("I112").Value = 3048.75
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 9146.25
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -8038.25
$ws.Range("N112").Value = $null

$ws.Range("H113").Value = 5940.174
$ws.Range("I113").Value = 687.3333
$ws.Range("J113").Value = 7794.1177
$ws.Range("K113").Value = 2061.9999
$ws.Range("L113").Value = 23382.3531
$ws.Range("M113").Value = 108.0001000000002
$ws.Range("N113").Value = -27722.3531

$ws.Range("H114").Value = 1882.619
$ws.Range("I114").Value = 1634.091
$ws.Range("J114").Value = 2156
$ws.Range("K114").Value = 4902.272999999999
$ws.Range("L114").Value = 6468
$ws.Range("M114").Value = -1648.272999999999
$ws.Range("N114").Value = -12976

$ws.Range("H116").Value = 1818665.4
$ws.Range("J116").Value = 4998
$ws.Range("L116").Value = 14994
$ws.Range("N116").Value = -21878

$ws.Range("H117").Value = 431.375
$ws.Range("I117").Value = 300.23077
$ws.Range("J117").Value = 999.6667
$ws.Range("K117").Value = 900.69231
$ws.Range("L117").Value = 2999.0001
$ws.Range("M117").Value = 2541.30769
$ws.Range("N117").Value = -9883.000100000001

$ws.Range("H119").Value = 50000
$ws.Range("I119").Value = 50000
$ws.Range("K119").Value = 150000
$ws.Range("M119").Value = -145162

$ws.Range("H120").Value = 13500
$ws.Range("I120").Value = 13500
$ws.Range("K120").Value = 40500
$ws.Range("M120").Value = -35662

$ws.Range("H123").Value = 1800
$ws.Range("J123").Value = 5800
$ws.Range("L123").Value = 17400
$ws.Range("N123").Value = -22300

$ws.Range("H129").Value = 1653.7391
$ws.Range("I129").Value = 1103.5
$ws.Range("J129").Value = 2509.6667
$ws.Range("K129").Value = 3310.5
$ws.Range("L129").Value = 7529.000100000001
$ws.Range("M129").Value = 1689.5
$ws.Range("N129").Value = -17529.0001

$ws.Range("H131").Value = 4611.3335
$ws.Range("I131").Value = 14617.875
$ws.Range("J131").Value = 2029
$ws.Range("K131").Value = 43853.625
$ws.Range("L131").Value = 6087
$ws.Range("M131").Value = -38813.625
$ws.Range("N131").Value = -16167

$ws.Range("H135").Value = 2789.8333
$ws.Range("I135").Value = 576.375
$ws.Range("J135").Value = 3594.7273
$ws.Range("K135").Value = 5187.375
$ws.Range("L135").Value = 32352.5457
$ws.Range("M135").Value = -2652.375
$ws.Range("N135").Value = -37422.5457

$ws.Range("H136").Value = 790
$ws.Range("I136").Value = 790
$ws.Range("K136").Value = 2370
$ws.Range("M136").Value = 2730

$ws.Range("H140").Value = 10082.839
$ws.Range("I140").Value = 4931.857
$ws.Range("J140").Value = 20899.9
$ws.Range("K140").Value = 14795.571
$ws.Range("L140").Value = 62699.7
$ws.Range("M140").Value = -9615.571
$ws.Range("N140").Value = -73059.70000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 368749.38
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 158332.5
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 158332.5
$ws.Range("M11").Value = -999861
$ws.Range("N11").Value = -158610.5

$ws.Range("H70").Value = 16198383
$ws.Range("I70").Value = 22822584
$ws.Range("K70").Value = 22822584
$ws.Range("M70").Value = -22822314

$ws.Range("H73").Value = 16198383
$ws.Range("I73").Value = 22822584
$ws.Range("K73").Value = 22822584
$ws.Range("M73").Value = -22821648

$ws.Range("H80").Value = 76925390
$ws.Range("I80").Value = 166668110
$ws.Range("J80").Value = 3056.5715
$ws.Range("K80").Value = 166668110
$ws.Range("L80").Value = 3056.5715
$ws.Range("M80").Value = -166667112
$ws.Range("N80").Value = -5052.5715

$ws.Range("H83").Value = 76925390
$ws.Range("I83").Value = 166668110
$ws.Range("J83").Value = 3056.5715
$ws.Range("K83").Value = 833340550
$ws.Range("L83").Value = 15282.8575
$ws.Range("M83").Value = -833335558
$ws.Range("N83").Value = -25266.8575

$ws.Range("H99").Value = 8453.4
$ws.Range("I99").Value = 1187
$ws.Range("J99").Value = 37519
$ws.Range("K99").Value = 1187
$ws.Range("L99").Value = 37519
$ws.Range("M99").Value = 1059
$ws.Range("N99").Value = -42011

$ws.Range("H122").Value = 7700301.5
$ws.Range("I122").Value = 15388613
$ws.Range("J122").Value = 11990
$ws.Range("K122").Value = 46165839
$ws.Range("L122").Value = 35970
$ws.Range("M122").Value = -46163389
$ws.Range("N122").Value = -40870

$ws.Range("H126").Value = 6602.8
$ws.Range("I126").Value = 2842.8
$ws.Range("J126").Value = 10362.8
$ws.Range("K126").Value = 8528.400000000001
$ws.Range("L126").Value = 31088.4
$ws.Range("M126").Value = -6058.400000000001
$ws.Range("N126").Value = -36028.39999999999

$ws.Range("H132").Value = 2059.175
$ws.Range("I132").Value = 1400.8438
$ws.Range("J132").Value = 4692.5
$ws.Range("K132").Value = 4202.5314
$ws.Range("L132").Value = 14077.5
$ws.Range("M132").Value = -1672.5314
$ws.Range("N132").Value = -19137.5

$ws.Range("H134").Value = 60016.75
$ws.Range("J134").Value = 60016.75
$ws.Range("L134").Value = 180050.25
$ws.Range("N134").Value = -185120.25

$ws.Range("H136").Value = 11124.353
$ws.Range("J136").Value = 11124.353
$ws.Range("L136").Value = 33373.05899999999
$ws.Range("N136").Value = -38473.05899999999

$ws.Range("H139").Value = 96214.44500000001
$ws.Range("J139").Value = 96214.44500000001
$ws.Range("L139").Value = 96214.44500000001
$ws.Range("N139").Value = -106494.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4581.5
$ws.Range("I7").Value = 3870.4348
$ws.Range("K7").Value = 3870.4348
$ws.Range("M7").Value = -3758.4348

$ws.Range("H22").Value = 611.25
$ws.Range("I22").Value = 577
$ws.Range("J22").Value = 714
$ws.Range("K22").Value = 577
$ws.Range("L22").Value = 714
$ws.Range("M22").Value = -282
$ws.Range("N22").Value = -1304

$ws.Range("H27").Value = 611.25
$ws.Range("I27").Value = 577
$ws.Range("J27").Value = 714
$ws.Range("K27").Value = 577
$ws.Range("L27").Value = 714
$ws.Range("M27").Value = -470
$ws.Range("N27").Value = -928

$ws.Range("H40").Value = 5284.3057
$ws.Range("I40").Value = 5457.3
$ws.Range("J40").Value = 4419.3335
$ws.Range("K40").Value = 5457.3
$ws.Range("L40").Value = 4419.3335
$ws.Range("M40").Value = -5321.3
$ws.Range("N40").Value = -4691.3335

$ws.Range("H46").Value = 1937.6666
$ws.Range("J46").Value = 1841.091
$ws.Range("L46").Value = 1841.091
$ws.Range("N46").Value = -2217.091

$ws.Range("H51").Value = 50084
$ws.Range("J51").Value = 50084
$ws.Range("L51").Value = 50084
$ws.Range("N51").Value = -51040

$ws.Range("H55").Value = 369.73685
$ws.Range("I55").Value = 283.73334
$ws.Range("K55").Value = 283.73334
$ws.Range("M55").Value = -110.73334

$ws.Range("H61").Value = 1644.3077
$ws.Range("I61").Value = 1526.762
$ws.Range("K61").Value = 1526.762
$ws.Range("M61").Value = -1324.762

$ws.Range("H68").Value = 9332.333000000001
$ws.Range("I68").Value = 9498.5
$ws.Range("J68").Value = 9000
$ws.Range("K68").Value = 9498.5
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = -8749.5
$ws.Range("N68").Value = -10498

$ws.Range("H71").Value = 9332.333000000001
$ws.Range("I71").Value = 9498.5
$ws.Range("J71").Value = 9000
$ws.Range("K71").Value = 47492.5
$ws.Range("L71").Value = 45000
$ws.Range("M71").Value = -43748.5
$ws.Range("N71").Value = -52488

$ws.Range("H93").Value = 1917
$ws.Range("I93").Value = 2018.8889
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 2018.8889
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -770.8888999999999
$ws.Range("N93").Value = -3496

$ws.Range("H113").Value = 1644.3077
$ws.Range("I113").Value = 1526.762
$ws.Range("K113").Value = 1526.762
$ws.Range("M113").Value = 643.2380000000001

$ws.Range("H122").Value = 3983
$ws.Range("I122").Value = 4399.75
$ws.Range("J122").Value = 3149.5
$ws.Range("K122").Value = 13199.25
$ws.Range("L122").Value = 9448.5
$ws.Range("M122").Value = -10749.25
$ws.Range("N122").Value = -14348.5

$ws.Range("H126").Value = 4581.5
$ws.Range("I126").Value = 3870.4348
$ws.Range("K126").Value = 11611.3044
$ws.Range("M126").Value = -9141.304400000001

$ws.Range("H132").Value = 8264.200000000001
$ws.Range("I132").Value = 6277.5
$ws.Range("K132").Value = 18832.5
$ws.Range("M132").Value = -16302.5

$ws.Range("H133").Value = 102322.336
$ws.Range("J133").Value = 102322.336
$ws.Range("L133").Value = 102322.336
$ws.Range("N133").Value = -107382.336

$ws.Range("H135").Value = 67000
$ws.Range("J135").Value = 67000
$ws.Range("L135").Value = 67000
$ws.Range("N135").Value = -77140

$ws.Range("H138").Value = 88123
$ws.Range("J138").Value = 88123
$ws.Range("L138").Value = 88123
$ws.Range("N138").Value = -98403

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 38462210
$ws.Range("I4").Value = 52631936
$ws.Range("J4").Value = 1530.2858
$ws.Range("K4").Value = 52631936
$ws.Range("L4").Value = 1530.2858
$ws.Range("M4").Value = -52631823
$ws.Range("N4").Value = -1756.2858

$ws.Range("H100").Value = 142858600
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459

$ws.Range("H107").Value = 748
$ws.Range("I107").Value = 784.3333
$ws.Range("K107").Value = 2352.9999
$ws.Range("M107").Value = -432.9998999999998

$ws.Range("H113").Value = 587.625
$ws.Range("I113").Value = 542.05554
$ws.Range("J113").Value = 724.3333
$ws.Range("K113").Value = 1626.16662
$ws.Range("L113").Value = 2172.9999
$ws.Range("M113").Value = 543.83338
$ws.Range("N113").Value = -6512.9999

$ws.Range("H132").Value = 15155562
$ws.Range("I132").Value = 16670869
$ws.Range("K132").Value = 50012607
$ws.Range("M132").Value = -50010077

$ws.Range("H135").Value = 54520.777
$ws.Range("J135").Value = 54520.777
$ws.Range("L135").Value = 54520.777
$ws.Range("N135").Value = -64660.777

$ws.Range("H140").Value = 129179
$ws.Range("J140").Value = 129179
$ws.Range("L140").Value = 129179
$ws.Range("N140").Value = -139539

$ws.Range("H141").Value = 69998.86
$ws.Range("J141").Value = 69998.86
$ws.Range("L141").Value = 69998.86
$ws.Range("N141").Value = -80358.86
